$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label
$ws.Range("E1").Value = "strength (RMS)"

# Update data values (columns B, D, E) for rows 2-19
$data = @{
    2  = @{ B = 16.4;  D = 13.2;  E = 51.8 }
    3  = @{ B = 14;    D = 13.6;  E = 62.2 }
    4  = @{ E = 62 }
    5  = @{ B = 16.8;  D = 13.6;  E = 49.8 }
    6  = @{ B = 13.6;  D = 14;    E = 67.40000000000001 }
    7  = @{ B = 14;    D = 14;    E = 53.5 }
    8  = @{ B = 14.8;  D = 19.6;  E = 56.6 }
    9  = @{ B = 14;    C = 35.5;  D = 21.5;  E = 66 }
    10 = @{ B = 15.2;  D = 15.2;  E = 54.4 }
    11 = @{ B = 14;    D = 16;    E = 66.8 }
    12 = @{ B = 13.6;  D = 16.4;  E = 67.2 }
    13 = @{ B = 14.4;  D = 15.6;  E = 61.4 }
    14 = @{ B = 14;    D = 16;    E = 67.8 }
    15 = @{ B = 13;    D = 17;    E = 68.75 }
    16 = @{ E = 63 }
    17 = @{ B = 14.4;  D = 15.6;  E = 59.8 }
    18 = @{ B = 14.5;  D = 16;    E = 60.25 }
    19 = @{ B = 17.6;  D = 12.4;  E = 52.4 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
